$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 10
$ws.Range("I10").Value = 10
$ws.Range("K10").Value = 10
$ws.Range("M10").Value = 283
$ws.Range("H20").Value = 6103.5
$ws.Range("I20").Value = 524.2
$ws.Range("J20").Value = 34000
$ws.Range("K20").Value = 524.2
$ws.Range("L20").Value = 34000
$ws.Range("M20").Value = -294.2
$ws.Range("N20").Value = -34460
$ws.Range("H28").Value = 356.8889
$ws.Range("J28").Value = 357.75
$ws.Range("L28").Value = 357.75
$ws.Range("N28").Value = -1327.75
$ws.Range("H35").Value = 6103.5
$ws.Range("I35").Value = 524.2
$ws.Range("J35").Value = 34000
$ws.Range("K35").Value = 524.2
$ws.Range("L35").Value = 34000
$ws.Range("M35").Value = -145.2
$ws.Range("N35").Value = -34758
$ws.Range("H125").Value = 406.07693
$ws.Range("I125").Value = 869.3333
$ws.Range("K125").Value = 7823.9997
$ws.Range("M125").Value = -5363.9997
$ws.Range("H129").Value = 2375.8147
$ws.Range("I129").Value = 4049.25
$ws.Range("J129").Value = 2084.7827
$ws.Range("K129").Value = 12147.75
$ws.Range("L129").Value = 6254.348100000001
$ws.Range("M129").Value = -7147.75
$ws.Range("N129").Value = -16254.3481

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3737.4849
$ws.Range("I61").Value = 2311
$ws.Range("J61").Value = 4121.5386
$ws.Range("K61").Value = 2311
$ws.Range("L61").Value = 4121.5386
$ws.Range("M61").Value = -2099
$ws.Range("N61").Value = -4545.5386
$ws.Range("H102").Value = 13501.947
$ws.Range("I102").Value = 2785.889
$ws.Range("J102").Value = 23146.4
$ws.Range("K102").Value = 2785.889
$ws.Range("L102").Value = 23146.4
$ws.Range("M102").Value = -1163.889
$ws.Range("N102").Value = -26390.4
$ws.Range("H132").Value = 20001824
$ws.Range("I132").Value = 38462772
$ws.Range("J132").Value = 2466.3333
$ws.Range("K132").Value = 115388316
$ws.Range("L132").Value = 7398.999899999999
$ws.Range("M132").Value = -115385786
$ws.Range("N132").Value = -12458.9999
$ws.Range("H136").Value = 3737.4849
$ws.Range("I136").Value = 2311
$ws.Range("J136").Value = 4121.5386
$ws.Range("K136").Value = 6933
$ws.Range("L136").Value = 12364.6158
$ws.Range("M136").Value = -4383
$ws.Range("N136").Value = -17464.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3865.5376
$ws.Range("I134").Value = 3055.7058
$ws.Range("J134").Value = 4084.0635
$ws.Range("K134").Value = 9167.117400000001
$ws.Range("L134").Value = 12252.1905
$ws.Range("M134").Value = -6632.117400000001
$ws.Range("N134").Value = -17322.1905

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 15000
$ws.Range("J8").Value = 15000
$ws.Range("L8").Value = 15000
$ws.Range("N8").Value = -15280
$ws.Range("H12").Value = 17073.2
$ws.Range("J12").Value = 20666.5
$ws.Range("L12").Value = 20666.5
$ws.Range("N12").Value = -21006.5
$ws.Range("H50").Value = 39672.25
$ws.Range("J50").Value = 39672.25
$ws.Range("L50").Value = 39672.25
$ws.Range("N50").Value = -40922.25
$ws.Range("H60").Value = 21104.791
$ws.Range("I60").Value = 5797.6665
$ws.Range("J60").Value = 23291.523
$ws.Range("K60").Value = 5797.6665
$ws.Range("L60").Value = 23291.523
$ws.Range("M60").Value = -5286.6665
$ws.Range("N60").Value = -24313.523
$ws.Range("H134").Value = 47154.324
$ws.Range("I134").Value = 1202.5652
$ws.Range("K134").Value = 3607.6956
$ws.Range("M134").Value = -1072.6956

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 83334616
$ws.Range("J62").Value = 100001340
$ws.Range("L62").Value = 300004020
$ws.Range("N62").Value = -300005392
$ws.Range("H65").Value = 83334616
$ws.Range("J65").Value = 100001340
$ws.Range("L65").Value = 900012060
$ws.Range("N65").Value = -900018924
$ws.Range("H70").Value = 4568
$ws.Range("I70").Value = 1704
$ws.Range("K70").Value = 5112
$ws.Range("M70").Value = -4797
$ws.Range("H73").Value = 4568
$ws.Range("I73").Value = 1704
$ws.Range("K73").Value = 5112
$ws.Range("M73").Value = -4020
$ws.Range("H75").Value = 1581.5
$ws.Range("I75").Value = 663
$ws.Range("K75").Value = 1989
$ws.Range("M75").Value = -991
$ws.Range("H78").Value = 1581.5
$ws.Range("I78").Value = 663
$ws.Range("K78").Value = 5967
$ws.Range("M78").Value = -975
$ws.Range("H98").Value = 599
$ws.Range("I98").Value = 467.57144
$ws.Range("J98").Value = 714
$ws.Range("K98").Value = 1402.71432
$ws.Range("L98").Value = 2142
$ws.Range("M98").Value = 95.28567999999996
$ws.Range("N98").Value = -5138
$ws.Range("H103").Value = 1371.7142
$ws.Range("I103").Value = 193.5
$ws.Range("J103").Value = 2942.6667
$ws.Range("K103").Value = 580.5
$ws.Range("L103").Value = 8828.000100000001
$ws.Range("M103").Value = 298.5
$ws.Range("N103").Value = -10586.0001
$ws.Range("H114").Value = 1234.8889
$ws.Range("I114").Value = 762.2857
$ws.Range("J114").Value = 1535.6364
$ws.Range("K114").Value = 2286.8571
$ws.Range("L114").Value = 4606.9092
$ws.Range("M114").Value = 967.1428999999998
$ws.Range("N114").Value = -11114.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 4097.5
$ws.Range("J19").Value = 6000
$ws.Range("L19").Value = 6000
$ws.Range("N19").Value = -6576
$ws.Range("H80").Value = 189155.73
$ws.Range("I80").Value = 281611.94
$ws.Range("J80").Value = 4243.3335
$ws.Range("K80").Value = 281611.94
$ws.Range("L80").Value = 4243.3335
$ws.Range("M80").Value = -280613.94
$ws.Range("N80").Value = -6239.3335
$ws.Range("H83").Value = 189155.73
$ws.Range("I83").Value = 281611.94
$ws.Range("J83").Value = 4243.3335
$ws.Range("K83").Value = 1408059.7
$ws.Range("L83").Value = 21216.6675
$ws.Range("M83").Value = -1403067.7
$ws.Range("N83").Value = -31200.6675
$ws.Range("H132").Value = 4465.8
$ws.Range("I132").Value = 2082.647
$ws.Range("J132").Value = 6716.5557
$ws.Range("K132").Value = 6247.941
$ws.Range("L132").Value = 20149.6671
$ws.Range("M132").Value = -3717.941
$ws.Range("N132").Value = -25209.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H22").Value = 1206.909
$ws.Range("I22").Value = 799.3333
$ws.Range("J22").Value = 1696
$ws.Range("K22").Value = 799.3333
$ws.Range("L22").Value = 1696
$ws.Range("M22").Value = -504.3333
$ws.Range("N22").Value = -2286
$ws.Range("H27").Value = 1206.909
$ws.Range("I27").Value = 799.3333
$ws.Range("J27").Value = 1696
$ws.Range("K27").Value = 799.3333
$ws.Range("L27").Value = 1696
$ws.Range("M27").Value = -692.3333
$ws.Range("N27").Value = -1910
$ws.Range("H46").Value = 11067.083
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 800
$ws.Range("M46").Value = -612
$ws.Range("H134").Value = 42262.176
$ws.Range("J134").Value = 42262.176
$ws.Range("L134").Value = 42262.176
$ws.Range("N134").Value = -52402.176

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 31000
$ws.Range("J32").Value = 31000
$ws.Range("L32").Value = 31000
$ws.Range("N32").Value = -31634
$ws.Range("H132").Value = 1642.25
$ws.Range("I132").Value = 871.5625
$ws.Range("J132").Value = 3183.625
$ws.Range("K132").Value = 2614.6875
$ws.Range("L132").Value = 9550.875
$ws.Range("M132").Value = -84.6875
$ws.Range("N132").Value = -14610.875
$ws.Range("H133").Value = 78862.25
$ws.Range("J133").Value = 78862.25
$ws.Range("L133").Value = 78862.25
$ws.Range("N133").Value = -88982.25
